$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PIR")
$ws.Range($ws.Cells.Item(133, 1), $ws.Cells.Item(146, 1)).NumberFormat = "@"
$ws.Cells.Item(133, 1).Value = "2026-01-28"
$ws.Cells.Item(133, 2).Value = "16:21:32"
$ws.Cells.Item(133, 3).Value = "16:00"
$ws.Cells.Item(133, 4).Value = "Bathroom"
$ws.Cells.Item(133, 5).Value = "No Motion"
$ws.Cells.Item(133, 6).Value = "Inactive"
$ws.Cells.Item(134, 1).Value = "2026-01-28"
$ws.Cells.Item(134, 2).Value = "16:21:33"
$ws.Cells.Item(134, 3).Value = "16:00"
$ws.Cells.Item(134, 4).Value = "Bathroom"
$ws.Cells.Item(134, 5).Value = "No Motion"
$ws.Cells.Item(134, 6).Value = "Inactive"
$ws.Cells.Item(135, 1).Value = "2026-01-28"
$ws.Cells.Item(135, 2).Value = "16:21:36"
$ws.Cells.Item(135, 3).Value = "16:00"
$ws.Cells.Item(135, 4).Value = "Bathroom"
$ws.Cells.Item(135, 5).Value = "No Motion"
$ws.Cells.Item(135, 6).Value = "Inactive"
$ws.Cells.Item(136, 1).Value = "2026-01-28"
$ws.Cells.Item(136, 2).Value = "16:21:41"
$ws.Cells.Item(136, 3).Value = "16:00"
$ws.Cells.Item(136, 4).Value = "Bathroom"
$ws.Cells.Item(136, 5).Value = "No Motion"
$ws.Cells.Item(136, 6).Value = "Inactive"
$ws.Cells.Item(137, 1).Value = "2026-01-28"
$ws.Cells.Item(137, 2).Value = "16:21:46"
$ws.Cells.Item(137, 3).Value = "16:00"
$ws.Cells.Item(137, 4).Value = "Bathroom"
$ws.Cells.Item(137, 5).Value = "No Motion"
$ws.Cells.Item(137, 6).Value = "Inactive"
$ws.Cells.Item(138, 1).Value = "2026-01-28"
$ws.Cells.Item(138, 2).Value = "16:21:51"
$ws.Cells.Item(138, 3).Value = "16:00"
$ws.Cells.Item(138, 4).Value = "Bathroom"
$ws.Cells.Item(138, 5).Value = "No Motion"
$ws.Cells.Item(138, 6).Value = "Inactive"
$ws.Cells.Item(139, 1).Value = "2026-01-28"
$ws.Cells.Item(139, 2).Value = "16:21:56"
$ws.Cells.Item(139, 3).Value = "16:00"
$ws.Cells.Item(139, 4).Value = "Bathroom"
$ws.Cells.Item(139, 5).Value = "No Motion"
$ws.Cells.Item(139, 6).Value = "Inactive"
$ws.Cells.Item(140, 1).Value = "2026-01-28"
$ws.Cells.Item(140, 2).Value = "16:22:01"
$ws.Cells.Item(140, 3).Value = "16:00"
$ws.Cells.Item(140, 4).Value = "Bathroom"
$ws.Cells.Item(140, 5).Value = "No Motion"
$ws.Cells.Item(140, 6).Value = "Inactive"
$ws.Cells.Item(141, 1).Value = "2026-01-28"
$ws.Cells.Item(141, 2).Value = "16:22:06"
$ws.Cells.Item(141, 3).Value = "16:00"
$ws.Cells.Item(141, 4).Value = "Bathroom"
$ws.Cells.Item(141, 5).Value = "No Motion"
$ws.Cells.Item(141, 6).Value = "Inactive"
$ws.Cells.Item(142, 1).Value = "2026-01-28"
$ws.Cells.Item(142, 2).Value = "16:22:11"
$ws.Cells.Item(142, 3).Value = "16:00"
$ws.Cells.Item(142, 4).Value = "Bathroom"
$ws.Cells.Item(142, 5).Value = "No Motion"
$ws.Cells.Item(142, 6).Value = "Inactive"
$ws.Cells.Item(143, 1).Value = "2026-01-28"
$ws.Cells.Item(143, 2).Value = "16:22:16"
$ws.Cells.Item(143, 3).Value = "16:00"
$ws.Cells.Item(143, 4).Value = "Bathroom"
$ws.Cells.Item(143, 5).Value = "No Motion"
$ws.Cells.Item(143, 6).Value = "Inactive"
$ws.Cells.Item(144, 1).Value = "2026-01-28"
$ws.Cells.Item(144, 2).Value = "16:22:22"
$ws.Cells.Item(144, 3).Value = "16:00"
$ws.Cells.Item(144, 4).Value = "Bathroom"
$ws.Cells.Item(144, 5).Value = "No Motion"
$ws.Cells.Item(144, 6).Value = "Inactive"
$ws.Cells.Item(145, 1).Value = "2026-01-28"
$ws.Cells.Item(145, 2).Value = "16:22:26"
$ws.Cells.Item(145, 3).Value = "16:00"
$ws.Cells.Item(145, 4).Value = "Bathroom"
$ws.Cells.Item(145, 5).Value = "No Motion"
$ws.Cells.Item(145, 6).Value = "Inactive"
$ws.Cells.Item(146, 1).Value = "2026-01-28"
$ws.Cells.Item(146, 2).Value = "16:22:32"
$ws.Cells.Item(146, 3).Value = "16:00"
$ws.Cells.Item(146, 4).Value = "Bathroom"
$ws.Cells.Item(146, 5).Value = "No Motion"
$ws.Cells.Item(146, 6).Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range($ws.Cells.Item(131, 1), $ws.Cells.Item(143, 1)).NumberFormat = "@"
$ws.Range($ws.Cells.Item(131, 5), $ws.Cells.Item(143, 5)).NumberFormat = "@"
$ws.Cells.Item(131, 1).Value = "2026-01-28"
$ws.Cells.Item(131, 2).Value = "16:21:32"
$ws.Cells.Item(131, 3).Value = "16:00"
$ws.Cells.Item(131, 4).Value = "Bathroom"
$ws.Cells.Item(131, 5).Value = "87.4%"
$ws.Cells.Item(131, 6).Value = "Active"
$ws.Cells.Item(132, 1).Value = "2026-01-28"
$ws.Cells.Item(132, 2).Value = "16:21:33"
$ws.Cells.Item(132, 3).Value = "16:00"
$ws.Cells.Item(132, 4).Value = "Bathroom"
$ws.Cells.Item(132, 5).Value = "88.2%"
$ws.Cells.Item(132, 6).Value = "Active"
$ws.Cells.Item(133, 1).Value = "2026-01-28"
$ws.Cells.Item(133, 2).Value = "16:21:37"
$ws.Cells.Item(133, 3).Value = "16:00"
$ws.Cells.Item(133, 4).Value = "Bathroom"
$ws.Cells.Item(133, 5).Value = "88.3%"
$ws.Cells.Item(133, 6).Value = "Active"
$ws.Cells.Item(134, 1).Value = "2026-01-28"
$ws.Cells.Item(134, 2).Value = "16:21:41"
$ws.Cells.Item(134, 3).Value = "16:00"
$ws.Cells.Item(134, 4).Value = "Bathroom"
$ws.Cells.Item(134, 5).Value = "87.3%"
$ws.Cells.Item(134, 6).Value = "Active"
$ws.Cells.Item(135, 1).Value = "2026-01-28"
$ws.Cells.Item(135, 2).Value = "16:21:45"
$ws.Cells.Item(135, 3).Value = "16:00"
$ws.Cells.Item(135, 4).Value = "Bathroom"
$ws.Cells.Item(135, 5).Value = "88.2%"
$ws.Cells.Item(135, 6).Value = "Active"
$ws.Cells.Item(136, 1).Value = "2026-01-28"
$ws.Cells.Item(136, 2).Value = "16:21:53"
$ws.Cells.Item(136, 3).Value = "16:00"
$ws.Cells.Item(136, 4).Value = "Bathroom"
$ws.Cells.Item(136, 5).Value = "87.2%"
$ws.Cells.Item(136, 6).Value = "Active"
$ws.Cells.Item(137, 1).Value = "2026-01-28"
$ws.Cells.Item(137, 2).Value = "16:21:57"
$ws.Cells.Item(137, 3).Value = "16:00"
$ws.Cells.Item(137, 4).Value = "Bathroom"
$ws.Cells.Item(137, 5).Value = "88.2%"
$ws.Cells.Item(137, 6).Value = "Active"
$ws.Cells.Item(138, 1).Value = "2026-01-28"
$ws.Cells.Item(138, 2).Value = "16:22:01"
$ws.Cells.Item(138, 3).Value = "16:00"
$ws.Cells.Item(138, 4).Value = "Bathroom"
$ws.Cells.Item(138, 5).Value = "87.3%"
$ws.Cells.Item(138, 6).Value = "Active"
$ws.Cells.Item(139, 1).Value = "2026-01-28"
$ws.Cells.Item(139, 2).Value = "16:22:05"
$ws.Cells.Item(139, 3).Value = "16:00"
$ws.Cells.Item(139, 4).Value = "Bathroom"
$ws.Cells.Item(139, 5).Value = "88.2%"
$ws.Cells.Item(139, 6).Value = "Active"
$ws.Cells.Item(140, 1).Value = "2026-01-28"
$ws.Cells.Item(140, 2).Value = "16:22:13"
$ws.Cells.Item(140, 3).Value = "16:00"
$ws.Cells.Item(140, 4).Value = "Bathroom"
$ws.Cells.Item(140, 5).Value = "88.2%"
$ws.Cells.Item(140, 6).Value = "Active"
$ws.Cells.Item(141, 1).Value = "2026-01-28"
$ws.Cells.Item(141, 2).Value = "16:22:17"
$ws.Cells.Item(141, 3).Value = "16:00"
$ws.Cells.Item(141, 4).Value = "Bathroom"
$ws.Cells.Item(141, 5).Value = "88.2%"
$ws.Cells.Item(141, 6).Value = "Active"
$ws.Cells.Item(142, 1).Value = "2026-01-28"
$ws.Cells.Item(142, 2).Value = "16:22:21"
$ws.Cells.Item(142, 3).Value = "16:00"
$ws.Cells.Item(142, 4).Value = "Bathroom"
$ws.Cells.Item(142, 5).Value = "87.3%"
$ws.Cells.Item(142, 6).Value = "Active"
$ws.Cells.Item(143, 1).Value = "2026-01-28"
$ws.Cells.Item(143, 2).Value = "16:22:29"
$ws.Cells.Item(143, 3).Value = "16:00"
$ws.Cells.Item(143, 4).Value = "Bathroom"
$ws.Cells.Item(143, 5).Value = "88.2%"
$ws.Cells.Item(143, 6).Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Range($ws.Cells.Item(131, 1), $ws.Cells.Item(143, 1)).NumberFormat = "@"
$ws.Cells.Item(131, 1).Value = "2026-01-28"
$ws.Cells.Item(131, 2).Value = "16:21:33"
$ws.Cells.Item(131, 3).Value = "16:00"
$ws.Cells.Item(131, 4).Value = "Bathroom"
$ws.Cells.Item(131, 5).Value = "22.8C"
$ws.Cells.Item(131, 6).Value = "Active"
$ws.Cells.Item(132, 1).Value = "2026-01-28"
$ws.Cells.Item(132, 2).Value = "16:21:34"
$ws.Cells.Item(132, 3).Value = "16:00"
$ws.Cells.Item(132, 4).Value = "Bathroom"
$ws.Cells.Item(132, 5).Value = "22.7C"
$ws.Cells.Item(132, 6).Value = "Active"
$ws.Cells.Item(133, 1).Value = "2026-01-28"
$ws.Cells.Item(133, 2).Value = "16:21:37"
$ws.Cells.Item(133, 3).Value = "16:00"
$ws.Cells.Item(133, 4).Value = "Bathroom"
$ws.Cells.Item(133, 5).Value = "22.8C"
$ws.Cells.Item(133, 6).Value = "Active"
$ws.Cells.Item(134, 1).Value = "2026-01-28"
$ws.Cells.Item(134, 2).Value = "16:21:41"
$ws.Cells.Item(134, 3).Value = "16:00"
$ws.Cells.Item(134, 4).Value = "Bathroom"
$ws.Cells.Item(134, 5).Value = "22.8C"
$ws.Cells.Item(134, 6).Value = "Active"
$ws.Cells.Item(135, 1).Value = "2026-01-28"
$ws.Cells.Item(135, 2).Value = "16:21:45"
$ws.Cells.Item(135, 3).Value = "16:00"
$ws.Cells.Item(135, 4).Value = "Bathroom"
$ws.Cells.Item(135, 5).Value = "22.8C"
$ws.Cells.Item(135, 6).Value = "Active"
$ws.Cells.Item(136, 1).Value = "2026-01-28"
$ws.Cells.Item(136, 2).Value = "16:21:53"
$ws.Cells.Item(136, 3).Value = "16:00"
$ws.Cells.Item(136, 4).Value = "Bathroom"
$ws.Cells.Item(136, 5).Value = "22.7C"
$ws.Cells.Item(136, 6).Value = "Active"
$ws.Cells.Item(137, 1).Value = "2026-01-28"
$ws.Cells.Item(137, 2).Value = "16:21:57"
$ws.Cells.Item(137, 3).Value = "16:00"
$ws.Cells.Item(137, 4).Value = "Bathroom"
$ws.Cells.Item(137, 5).Value = "22.8C"
$ws.Cells.Item(137, 6).Value = "Active"
$ws.Cells.Item(138, 1).Value = "2026-01-28"
$ws.Cells.Item(138, 2).Value = "16:22:01"
$ws.Cells.Item(138, 3).Value = "16:00"
$ws.Cells.Item(138, 4).Value = "Bathroom"
$ws.Cells.Item(138, 5).Value = "22.8C"
$ws.Cells.Item(138, 6).Value = "Active"
$ws.Cells.Item(139, 1).Value = "2026-01-28"
$ws.Cells.Item(139, 2).Value = "16:22:05"
$ws.Cells.Item(139, 3).Value = "16:00"
$ws.Cells.Item(139, 4).Value = "Bathroom"
$ws.Cells.Item(139, 5).Value = "22.8C"
$ws.Cells.Item(139, 6).Value = "Active"
$ws.Cells.Item(140, 1).Value = "2026-01-28"
$ws.Cells.Item(140, 2).Value = "16:22:13"
$ws.Cells.Item(140, 3).Value = "16:00"
$ws.Cells.Item(140, 4).Value = "Bathroom"
$ws.Cells.Item(140, 5).Value = "22.8C"
$ws.Cells.Item(140, 6).Value = "Active"
$ws.Cells.Item(141, 1).Value = "2026-01-28"
$ws.Cells.Item(141, 2).Value = "16:22:17"
$ws.Cells.Item(141, 3).Value = "16:00"
$ws.Cells.Item(141, 4).Value = "Bathroom"
$ws.Cells.Item(141, 5).Value = "22.8C"
$ws.Cells.Item(141, 6).Value = "Active"
$ws.Cells.Item(142, 1).Value = "2026-01-28"
$ws.Cells.Item(142, 2).Value = "16:22:21"
$ws.Cells.Item(142, 3).Value = "16:00"
$ws.Cells.Item(142, 4).Value = "Bathroom"
$ws.Cells.Item(142, 5).Value = "22.8C"
$ws.Cells.Item(142, 6).Value = "Active"
$ws.Cells.Item(143, 1).Value = "2026-01-28"
$ws.Cells.Item(143, 2).Value = "16:22:29"
$ws.Cells.Item(143, 3).Value = "16:00"
$ws.Cells.Item(143, 4).Value = "Bathroom"
$ws.Cells.Item(143, 5).Value = "22.8C"
$ws.Cells.Item(143, 6).Value = "Active"
